$wb = $excel.ActiveWorkbook

# --- Rename second sheet ---
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "More Analyses"

# --- Sheet1 ("Sheet1"): selection moves, tab no longer selected ---
$ws1.Range("A24:M24").Select()

# --- Sheet2 ("More Analyses"): rebuild the data with new layout ---
# Row 1 headers (reordered + renamed)
$ws2.Range("B1").Value = "mother"
$ws2.Range("C1").Value = "mother %"
$ws2.Range("D1").Value = "father"
$ws2.Range("E1").Value = "father %"
$ws2.Range("F1").Value = "self"
$ws2.Range("G1").Value = "self %"

# Row 2: Non-sci
$ws2.Range("A2").Value = "Non-sci"
$ws2.Range("B2").Value = 498.95
$ws2.Range("C2").Value = 87.87
$ws2.Range("D2").Value = 497.66
$ws2.Range("E2").Value = 90.34
$ws2.Range("F2").Value = 487.76
$ws2.Range("G2").Value = 59.89

# Row 3: Sci, engineer/tech (was "Sci" row)
$ws2.Range("A3").Value = "Sci, engineer/tech"
$ws2.Range("B3").Value = 570.16999999999996
$ws2.Range("C3").Value = 1.38
$ws2.Range("D3").Value = 554.24
$ws2.Range("E3").Value = 7.3
$ws2.Range("F3").Value = 543.75
$ws2.Range("G3").Value = 16.63

# Row 4: Sci, medical (new row)
$ws2.Range("A4").Value = "Sci, medical"
$ws2.Range("B4").Value = 506.39
$ws2.Range("C4").Value = 10.75
$ws2.Range("D4").Value = 539.75
$ws2.Range("E4").Value = 2.36
$ws2.Range("F4").Value = 296.83
$ws2.Range("G4").Value = 23.48

# Row 6: second table headers - urban / not urban school community
$ws2.Range("A6").Value = "Self"
$ws2.Range("B6").Value = "not urban"
$ws2.Range("C6").Value = "urban"
$ws2.Range("D6").Value = "not-urban %"
$ws2.Range("E6").Value = "urban %"

# Row 7: Non-sci
$ws2.Range("A7").Value = "Non-sci"
$ws2.Range("B7").Value = 497.45
$ws2.Range("C7").Value = 472.61
$ws2.Range("D7").Value = 60.98
$ws2.Range("E7").Value = 39.020000000000003

# Row 8: Sci, engineer/tech
$ws2.Range("A8").Value = "Sci, engineer/tech"
$ws2.Range("B8").Value = 552.86
$ws2.Range("C8").Value = 530.66999999999996
$ws2.Range("D8").Value = 58.95
$ws2.Range("E8").Value = 41.05

# Row 9: Sci, medical
$ws2.Range("A9").Value = "Sci, medical"
$ws2.Range("B9").Value = 503.99
$ws2.Range("C9").Value = 484.06
$ws2.Range("D9").Value = 64.08
$ws2.Range("E9").Value = 35.92

$ws2.Range("F13").Select()

# --- Workbook-level: activate "More Analyses" as the visible/active tab ---
$ws2.Activate()
